$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH124"
$ws.Range("C2").Value = "CONFLICTING NOTIONS OF SOVEREIGNTY, INDEPENDANT ELECTORAL COMMISSION DOCUMENTS AND ARMBAND, GOVERNMENT GAZETTE"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

$ws.Range("A2:H2").Font.Name = "Calibri"
